$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Use case1: cell A8 now holds the text "I" (was the number 7), and it
# becomes the new active selection on the sheet.
$ws.Range("A8").Value = "I"
$ws.Range("A8").Select()

# Page setup was touched as well (orientation explicitly set to portrait).
$ws.PageSetup.Orientation = 1
